$d = $word.ActiveDocument
$startPara = 24
$endPara = $d.Paragraphs.Count
$start = $d.Paragraphs.Item($startPara).Range.Start
$end = $d.Paragraphs.Item($endPara).Range.End
$r = $d.Range($start, $end)
$xml = '<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main'' w:rsidR=''00EF32BD'' w:rsidRPr=''00CC4DF5'' w:rsidRDefault=''00EF32BD'' w:rsidP=''00EF32BD''><w:pPr><w:ind w:firstLine=''708''/></w:pPr><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:proofErr w:type=''spellStart''/><w:proofErr w:type=''gramStart''/><w:r><w:t>transform.SetPositionAndRotation</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=''gramEnd''/><w:r><w:t>posición, rotación)</w:t></w:r></w:p><w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:ind w:left=''1413''/></w:pPr></w:p><w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:ind w:left=''1413''/></w:pPr><w:proofErr w:type=''spellStart''/><w:proofErr w:type=''gramStart''/><w:r><w:lastRenderedPageBreak/><w:t>Translate</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=''gramEnd''/><w:r><w:t>Vector): Mueve el Objeto en la dirección y distancia que apunta el Vector.</w:t></w:r></w:p><w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:ind w:left=''1413''/></w:pPr><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:proofErr w:type=''spellStart''/><w:proofErr w:type=''gramStart''/><w:r><w:t>transform.Translate</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=''gramEnd''/><w:r><w:t>Vector3.forward)</w:t></w:r></w:p><w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:ind w:left=''1413''/></w:pPr><w:proofErr w:type=''spellStart''/><w:proofErr w:type=''gramStart''/><w:r><w:t>Rotate</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=''gramEnd''/><w:r><w:t>Vector</w:t></w:r><w:r><w:t xml:space=''preserve''>, </w:t></w:r><w:proofErr w:type=''spellStart''/><w:r><w:t>float</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:t xml:space=''preserve''>, </w:t></w:r><w:proofErr w:type=''spellStart''/><w:r><w:t>GameObject</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:t>):</w:t></w:r><w:r><w:t xml:space=''preserve''> Ro</w:t></w:r><w:r><w:t xml:space=''preserve''>ta un objeto en el </w:t></w:r><w:proofErr w:type=''spellStart''/><w:r><w:t>angulo</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:t xml:space=''preserve''> que se le pase en el Vector, los grados que se le pasen en el </w:t></w:r><w:proofErr w:type=''spellStart''/><w:r><w:t>float</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:t xml:space=''preserve''> y sobre el </w:t></w:r><w:proofErr w:type=''spellStart''/><w:r><w:t>gameObject</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:t xml:space=''preserve''> relativo. Este </w:t></w:r><w:proofErr w:type=''spellStart''/><w:proofErr w:type=''gramStart''/><w:r><w:t>ultimo</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:proofErr w:type=''gramEnd''/><w:r><w:t xml:space=''preserve''> si no se pasa, por defecto es el objeto mismo.</w:t></w:r></w:p><w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:ind w:left=''2121'' w:firstLine=''3''/><w:rPr><w:lang w:val=''en-US''/></w:rPr></w:pPr><w:proofErr w:type=''spellStart''/><w:r><w:rPr><w:lang w:val=''en-US''/></w:rPr><w:t>transform.Rotate</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:rPr><w:lang w:val=''en-US''/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type=''spellStart''/><w:r><w:rPr><w:lang w:val=''en-US''/></w:rPr><w:t>Vector.up</w:t></w:r><w:proofErr w:type=''spellEnd''/><w:r><w:rPr><w:lang w:val=''en-US''/></w:rPr><w:t>, 5.0f</w:t></w:r><w:r><w:rPr><w:lang w:val=''en-US''/></w:rPr><w:t>)</w:t></w:r><w:bookmarkStart w:id=''0'' w:name=''_GoBack''/><w:bookmarkEnd w:id=''0''/></w:p><w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:rPr><w:lang w:val=''en-US''/></w:rPr></w:pPr></w:p><w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:rPr><w:lang w:val=''en-US''/></w:rPr></w:pPr></w:p>'
$result = $r.InsertXML($xml)
Write-Host "InsertXML result:" $result
